# Applies the "Saldo" export update:
#  - Insert a new account row for ALPHASITIO above BRASFORT
#  - Replace the ANUAR / JOAQUIM rows with GABRIEL / MARIA (positive balance) rows
#  - Insert a new account row for PEDRO above ANA
#  - Remove the trailing BRUNO and MARIA (negative balance) rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-AccountRow($Row, $Conta, $Nome, $Saldo) {
    $acctCell = $ws.Cells.Item($Row, 1)
    $acctCell.NumberFormat = "@"
    $acctCell.Value = $Conta
    $ws.Cells.Item($Row, 2).Value = $Nome
    $ws.Cells.Item($Row, 3).Value = $Saldo
}

# 1) Insert ALPHASITIO row above row 3 (BRASFORT)
$ws.Rows.Item(3).Insert()
Set-AccountRow 3 "005305448" "ALPHASITIO" 321659.51

# 2) Replace ANUAR (now row 5) with GABRIEL
Set-AccountRow 5 "005666419" "GABRIEL" 125000

# 3) Replace JOAQUIM (now row 6) with MARIA (positive balance)
Set-AccountRow 6 "004212581" "MARIA" 44449.83

# 4) Insert PEDRO row above row 14 (ANA, after the prior insert shifted it down)
$ws.Rows.Item(14).Insert()
Set-AccountRow 14 "005324840" "PEDRO" 4000

# 5) Remove the trailing BRUNO and MARIA (negative balance) rows (now rows 164-165)
$ws.Rows.Item(164).Delete()
$ws.Rows.Item(164).Delete()
